# Regenerate orders with updated distance/size labels.
# Mapping: D64->D69, D80->D86, D51->D55, S30->S31
# Applied to Condition, Filename_Left, Filename_Right, Distance, Size columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

function Convert-Label([string]$text) {
    if ($null -eq $text) { return $text }
    $result = $text
    $result = $result.Replace("D64", "`u{0001}")
    $result = $result.Replace("D80", "`u{0002}")
    $result = $result.Replace("D51", "`u{0003}")
    $result = $result.Replace("`u{0001}", "D69")
    $result = $result.Replace("`u{0002}", "D86")
    $result = $result.Replace("`u{0003}", "D55")
    $result = $result.Replace("S30", "S31")
    return $result
}

# Columns to update: B=Condition, D=Filename_Left, E=Filename_Right, H=Distance, J=Size
$columns = @(2, 4, 5, 8, 10)

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $columns) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = Convert-Label $val
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
